# A new daily price-report row for "Mango" at "Vega Modelo de Temuco" is
# inserted at row 458, pushing the existing rows 458:539 down to 459:540
# (sheet dimension grows from A1:T539 to A1:T540).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 458, shifting rows 458-539 down.
$ws.Rows(458).Insert()

# Populate the newly inserted row 458 with the new observation.
$ws.Range("A458").Value = 10
$ws.Range("B458").Value = "Vega Modelo de Temuco"
$ws.Range("C458").Value = "La Araucanía"
$ws.Range("D458").Value = 45015
$ws.Range("E458").Value = 9
$ws.Range("F458").Value = "Fruta"
$ws.Range("G458").Value = 100108
$ws.Range("H458").Value = "Tropicales y subtropicales"
$ws.Range("I458").Value = 100108002
$ws.Range("J458").Value = "Mango"
$ws.Range("K458").Value = "Sin especificar"
$ws.Range("L458").Value = "Primera"
$ws.Range("M458").Value = 530
$ws.Range("N458").Value = 7000
$ws.Range("O458").Value = 8000
$ws.Range("P458").Value = 7472
$ws.Range("Q458").Value = "`$/bandeja 4 kilos"
$ws.Range("R458").Value = "Perú"
$ws.Range("S458").Value = 1868
$ws.Range("T458").Value = 4
